# Generate Report for Handoff
# Refresh the "Latest Handoff Date/Datetime" timestamps for the files that
# were re-handed-off, on the Overview sheet (column D) and on each
# per-language sheet (zh-cn / de-de, column E).

$wb = $excel.ActiveWorkbook

# File names (in row order) whose handoff timestamp was refreshed.
$files = @(
    "a19684ea-365c-4503-950b-f2e47b4de564.md",
    "1b9746d4-e458-4318-a9cc-a3d06d728cbd.md",
    "4f02c3c1-e133-4b91-8fd8-132c2b173435.md",
    "63bdfb6f-e2a9-4907-bcc8-5fb8cdf36c9d.md",
    "64e00640-fa88-43f4-8e33-d00beed7a25e.md",
    "844c9c99-7b0a-4476-ba06-c62299c08379.md",
    "a41b6138-fe38-4f5f-af59-8908c5099d95.md",
    "d44fcb6e-9b32-40a8-b504-7e89e1dfe160.md"
)

# Overview sheet: column D = "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
for ($r = 2; $r -le 16; $r++) {
    $name = $wsOverview.Cells.Item($r, 1).Value()
    if ($files -contains $name) {
        $wsOverview.Cells.Item($r, 4).Value = "2016-03-24 02:31:59"
    }
}

# zh-cn sheet: column E = "Latest Handoff Datetime"
$wsZh = $wb.Worksheets.Item("zh-cn")
for ($r = 2; $r -le 16; $r++) {
    $name = $wsZh.Cells.Item($r, 1).Value()
    if ($files -contains $name) {
        $wsZh.Cells.Item($r, 5).Value = "2016-03-24 02:31:55"
    }
}

# de-de sheet: column E = "Latest Handoff Datetime"
$wsDe = $wb.Worksheets.Item("de-de")
for ($r = 2; $r -le 16; $r++) {
    $name = $wsDe.Cells.Item($r, 1).Value()
    if ($files -contains $name) {
        $wsDe.Cells.Item($r, 5).Value = "2016-03-24 02:31:59"
    }
}
